# DRP Cases - add analysis/recommendation notes + hyperlinks to the
# "Social networks" and "Mobile applications" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Social networks
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Social networks")

# Hyperlink the two Facebook-group URLs that were actually re-checked
# (they keep pointing at their own cell text).
$ws2.Hyperlinks.Add($ws2.Range("B5"), $ws2.Range("B5").Value2)
$ws2.Hyperlinks.Add($ws2.Range("B8"), $ws2.Range("B8").Value2)

# Header cell F1 ("Analysis") becomes wrap-text.
$ws2.Range("F1").WrapText = $true

# Per-row analysis text (column F) + recommendation (column G).
$ws2.Range("F2").Value = "Personal account of user on Facebook with no suspicious activity observed"
$ws2.Range("F3").Value = "After checking it was a group related to job offerings and hiring with job posting on ADIB"
$ws2.Range("F4").Value = "After checking it was a group related to job offerings and hiring with job posting on ADIB"
$ws2.Range("F5").Value = "After checking it was a group related to job offerings and hiring with job posting on ADIB"
$ws2.Range("F6").Value = "After checking it was found a group related to ADIB UAE with videos on offered services and nothing malicious was found"
$ws2.Range("F7").Value = "Personal account of user on Facebook with no suspicious activity observed "
$ws2.Range("F8").Value = "Group created on Facebook related to IT teams and nothing suspicious was found"
$ws2.Range("F2:F8").WrapText = $true

$ws2.Range("G9").Copy() | Out-Null
$ws2.Range("G2:G8").PasteSpecial(-4122) | Out-Null
$ws2.Range("G2:G8").Value = "Reject"
$excel.CutCopyMode = 0

# Row heights (wrapped text rows grow taller).
$ws2.Rows.Item(2).RowHeight = 72
$ws2.Rows.Item(3).RowHeight = 86.4
$ws2.Rows.Item(4).RowHeight = 86.4
$ws2.Rows.Item(5).RowHeight = 86.4
$ws2.Rows.Item(6).RowHeight = 115.2
$ws2.Rows.Item(7).RowHeight = 72
$ws2.Rows.Item(8).RowHeight = 72

# Move the active selection to G2 (where the reviewer ended up).
$ws2.Range("G2").Select()

# ---------------------------------------------------------------------
# Sheet: Mobile applications
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Mobile applications")

# Row 2: just an analysis note (a pasted hash value, its own little
# formatting - small Roboto font, no border/number-format, wrapped).
$c = $ws3.Range("F2")
$c.Value = "ff79f017f7c3d7d4746f736477780cf3b7e55eb0ae07cb56440787acd9985709"
$c.Style = "Normal"
$c.Font.Name = "Roboto"
$c.Font.Size = 8
$c.Font.Color = 2892832
$c.WrapText = $true
$ws3.Rows.Item(2).RowHeight = 21.6

# Rows 3-9: analysis (F) + recommendation (G) text. The whole B:F block
# on these rows also picks up an (empty) alignment record because the
# reviewer wrapped & then un-wrapped the text while formatting them.
$ws3.Range("B3:F9").WrapText = $true
$ws3.Range("B3:F9").WrapText = $false

$ws3.Range("F3").Value = "After checking link it was found broken to download ADIB Pay APP"
$ws3.Range("F4").Value = "After checking website link to download the app was not working and no suspicious activity was detected."
$ws3.Range("F5").Value = "After checking website link to download the app was not working and no suspicious activity was detected."
$ws3.Range("F6").Value = "After checking website link to download the app was not working and no suspicious activity was detected."
$ws3.Range("F7").Value = "After checking website link to download the app was not working and no suspicious activity was detected."
$ws3.Range("F8").Value = "After checking website link to download the app was not working and no suspicious activity was detected."
$ws3.Range("F9").Value = "Link was found not working with no suspicious actvity found"

$ws3.Range("G3:G9").Copy() | Out-Null
$ws3.Range("G3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws3.Range("G3").Value = "Reject"
$ws3.Range("G3").Copy() | Out-Null
$ws3.Range("G4:G9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws3.Range("G4:G9").Value = "Reject"

# The dubapk.com download link (row 8) gets hyperlinked to itself too.
$ws3.Hyperlinks.Add($ws3.Range("B8"), $ws3.Range("B8").Value2)

# Reviewer's final view state on this (active) sheet.
$ws3.Activate()
$ws3.Range("C10").Select()

$wb.Save()
